# Sprint 3 update for the SCRUM board: replace the Sprint 2 task cards with
# the new Sprint 3 tasks, moving every card into column A (single column of
# index cards instead of split across B/C), widening column A to fit, and
# moving the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any existing text from rows 3-13 across all three columns first.
foreach ($r in 3..13) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 3).ClearContents()
}

# New Sprint 3 task text, one per row (rows 3-11). Rows 12-13 stay blank.
# Order below matches the order the text was authored in (row 6's text was
# entered last), so new shared-string entries land in the same slots as the
# original edit.
$taskOrder = @(
    @{ Row = 3;  Text = "Generate a Venmo URL that can be used to charge everyone " }
    @{ Row = 4;  Text = "Use the amount of people to generate an editable form " }
    @{ Row = 5;  Text = "As a stopgap, use a student website that we are all provided with " }
    @{ Row = 7;  Text = "Update the cost divide function to calculate a variable cost" }
    @{ Row = 8;  Text = "Generate a Venmo URL for each user with variable cost " }
    @{ Row = 9;  Text = "Make sure appropriate values are entered in the initial form" }
    @{ Row = 10; Text = "In a custom split, make sure the custom split adds up to the total amount " }
    @{ Row = 11; Text = "Somehow attempt to verity that the Venmo usernames are legitimate" }
    @{ Row = 6;  Text = "Make the editable form of Venmo usernames also have a field for the %" }
)

foreach ($task in $taskOrder) {
    $cell = $ws.Cells.Item($task.Row, 1)
    $cell.Value = $task.Text
    # Cards used to wrap text inside a bordered box (column B/C style); now
    # they are plain, unbordered, non-wrapping cells in column A.
    $cell.WrapText = $false
    $cell.Borders.LineStyle = -4142
}

# Row heights were only tall because of wrapped text - restore auto height.
foreach ($r in 3..13) {
    $ws.Rows.Item($r).AutoFit()
}

# Column A needs to be widened to fit the longer task descriptions.
$ws.Columns.Item(1).ColumnWidth = 61.43

# Move the active selection to A23 (matches the saved selection state).
$ws.Range("A23").Select()

# Plain portrait orientation page setup.
$ps = $ws.PageSetup
$ps.Orientation = 1
$ps.FirstPageNumber = 1
